$d = $word.ActiveDocument

# Replace the placeholder ID text, and drop the trailing run that only
# contained a single space (the Find/Replace below consumes it so the
# paragraph ends up with just one run).
$d.Content.Find.Execute("**ID__AFFARS_5349_topic_2__ID** ", $false, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_SUBPART_5349_1__ID**", 2)

# Update the first paragraph's formatting: indent + paragraph border.
$p = $d.Paragraphs(1)
$p.Range.ParagraphFormat.LeftIndent = 11.25
$p.Range.ParagraphFormat.Borders.DistanceFromTop = 5
$p.Range.ParagraphFormat.Borders.DistanceFromLeft = 5
$p.Range.ParagraphFormat.Borders.DistanceFromBottom = 5
$p.Range.ParagraphFormat.Borders.DistanceFromRight = 5
